$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LocalStiffness")
$ws.Activate() | Out-Null

# Update header label: "freq(Hz)" -> "freq(rpm)"
$ws.Range("A2").Value = "freq(rpm)"

# Update frequency values (Hz -> rpm) in column A
$ws.Range("A4").Value = 500
$ws.Range("A5").Value = 750
$ws.Range("A6").Value = 1000
$ws.Range("A7").Value = 1250
$ws.Range("A8").Value = 1500

# Widen column A to fit the new label (match width of neighboring data columns)
$ws.Columns("A:A").ColumnWidth = 10.5

# Move the active selection to A9
$ws.Range("A9").Select() | Out-Null
